$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$grid = @(
    @(3,5,2,8,4,7,6,9,1),
    @(6,7,4,9,1,2,5,3,8),
    @(8,9,1,5,3,6,2,7,4),
    @(5,2,8,7,9,4,1,6,3),
    @(4,1,7,6,8,3,9,2,5),
    @(9,3,6,2,5,1,4,8,7),
    @(1,4,9,3,6,8,7,5,2),
    @(2,8,5,1,7,9,3,4,6),
    @(7,6,3,4,2,5,8,1,9)
)

for ($r = 0; $r -lt 9; $r++) {
    for ($c = 0; $c -lt 9; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $grid[$r][$c]
    }
}
